$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "283.40"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.37"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.413"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06206"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.591"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.565"
$ws.Range("E7").Value = "6KuCoinTokenKCS"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.508"
$ws.Range("E8").Value = "7FTXTokenFTT"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8236"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01399"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1658"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03470"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03233"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09196"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.735"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001661"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04761"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006515"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.006184"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001073"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001611"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.832"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.360"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3349"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1222"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04727"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007198"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004033"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1105"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01167"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006920"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000755"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.107"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003016"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00001410"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01248"
